# Update the cached "last updated" date shown by the auto date placeholder
# (fld type="datetimeFigureOut") on the slide master and on every slide
# layout, and reposition the logo group on the single slide.

$p = $ppt.ActivePresentation

$newDate = "10/01/2019"

# -- Slide master: date placeholder is shape index 3 --------------------
$master = $p.Slides.Item(1).Design.SlideMaster
$master.Shapes.Item(3).TextFrame.TextRange.Text = $newDate

# -- Every slide layout: date placeholder is also shape index 3 ---------
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    $layout.Shapes.Item(3).TextFrame.TextRange.Text = $newDate
}

# -- Slide 1: move the logo group (shape "Groupe 15") --------------------
$slide = $p.Slides.Item(1)
$logoGroup = $slide.Shapes.Item(1)
$logoGroup.Left = 332.5788188976378
$logoGroup.Top = 63.94465566929134
